$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 9: A9 was stored as a text string; convert it to a real
# Excel date/time serial number (matching the numeric date cells used
# elsewhere in column A), keeping the same displayed value.
$ws.Range("A9").Value = 44548.45486111111
$ws.Range("A9").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- New row 10 ---
$ws.Range("A10").Value = 44548.4875
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = "Buy"
$ws.Range("C10").Value = "AAPL"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 169.93
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 169.93
$ws.Range("H10").Value = -169.93
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = -680.4399999999998
$ws.Range("L10").Value = -510.5099999999998
$ws.Range("M10").Value = -170.1699999999999
$ws.Range("N10").Value = "'"
$ws.Range("N10").Style = "Normal"
$ws.Range("O10").Value = "'"
$ws.Range("O10").Style = "Normal"
$ws.Range("P10").Value = "'"
$ws.Range("P10").Style = "Normal"
$ws.Range("Q10").Value = "'"
$ws.Range("Q10").Style = "Normal"

# --- New row 11 ---
$ws.Range("A11").Value = 44548.5
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = "Buy"
$ws.Range("C11").Value = "AMZN"
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 3354.21
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 3354.21
$ws.Range("H11").Value = -3354.21
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 3354.21
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = "'"
$ws.Range("N11").Style = "Normal"
$ws.Range("O11").Value = "'"
$ws.Range("O11").Style = "Normal"
$ws.Range("P11").Value = "'"
$ws.Range("P11").Style = "Normal"
$ws.Range("Q11").Value = "'"
$ws.Range("Q11").Style = "Normal"

# --- New row 12 (date kept as plain text, like the original row 9) ---
$ws.Range("A12").Value = "2021-12-18 12:01:00"
$ws.Range("B12").Value = "Buy"
$ws.Range("C12").Value = "AMZN"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 3354.21
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 3354.21
$ws.Range("H12").Value = -3354.21
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 2
$ws.Range("K12").Value = 3354.21
$ws.Range("L12").Value = 6708.42
$ws.Range("M12").Value = 3354.21
$ws.Range("N12").Value = "'"
$ws.Range("N12").Style = "Normal"
$ws.Range("O12").Value = "'"
$ws.Range("O12").Style = "Normal"
$ws.Range("P12").Value = "'"
$ws.Range("P12").Style = "Normal"
$ws.Range("Q12").Value = "'"
$ws.Range("Q12").Style = "Normal"
